# Insert a new data row for 2026/02/02 08:00 just before the current
# "2026/12/29" row (row 767). This pushes the existing rows 767-808 down to
# 768-809, growing the used range from A1:D808 to A1:D809.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(767).Insert()

# Column A stores dates as literal text (e.g. "2026/02/02"), not real date
# serials. Force a text number format before writing so Excel doesn't
# auto-convert the string into a date value, then restore the default
# "Normal" style so the new row matches the unstyled data rows around it.
$ws.Range("A767").NumberFormat = "@"
$ws.Range("A767").Value = "2026/02/02"
$ws.Range("B767").Value = "月"
$ws.Range("C767").Value = 8
$ws.Range("D767").Value = 162
$ws.Range("A767").Style = "Normal"
